$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.841.15"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
$ws.Range("D3").Value = "2.230.10"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "274.11"
$ws.Range("E5").Value = "  +6.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.92"
$ws.Range("E6").Value = "  +9.53%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -2.00%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -0.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.03"
$ws.Range("E10").Value = "  +3.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  -2.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.57"
$ws.Range("E12").Value = "  +6.56%  "

# Row 13
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").Value = "2.563.16"
$ws.Range("E14").Value = "  -0.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.05"
$ws.Range("E15").Value = "  +1.95%  "

# Row 16
$ws.Range("D16").Value = "2.247.24"
$ws.Range("E16").Value = "  -1.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "43.742.21"
$ws.Range("E18").Value = "  -0.72%  "

# Row 19
$ws.Range("E19").Value = "  -1.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.17"
$ws.Range("E20").Value = "  -2.08%  "

# Row 21
$ws.Range("E21").Value = "  -1.94%  "

# Row 22
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.27"
$ws.Range("E23").Value = "  -1.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.68"
$ws.Range("E24").Value = "  -7.33%  "

# Row 25
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("E26").Value = "  +14.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -1.62%  "

# Row 28
$ws.Range("E28").Value = "  +3.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +3.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.09"
$ws.Range("E30").Value = "  -3.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.79"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0906"
$ws.Range("E32").Value = "  +3.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.74"
$ws.Range("E33").Value = "  +0.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("E36").Value = "  -2.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("E37").Value = "  -4.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.27"
$ws.Range("E38").Value = "  -5.59%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  +20.01%  "

# Row 40
$ws.Range("E40").Value = "  +0.48%  "

# Row 41
$ws.Range("E41").Value = "  -4.57%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.209"
$ws.Range("E42").Value = "  +2.57%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.61"
$ws.Range("E43").Value = "  +0.39%  "

# Row 44
$ws.Range("E44").Value = "  -2.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.47"
$ws.Range("E45").Value = "  -1.20%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0981"
$ws.Range("E46").Value = "  -0.91%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.16"
$ws.Range("E47").Value = "  -3.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  +0.70%  "

# Row 49
$ws.Range("E49").Value = "  +2.90%  "

# Row 50
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.431"
$ws.Range("E50").Value = "  -5.82%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.49"
$ws.Range("E51").Value = "  -2.20%  "
